# Update "想去人数" (F column) and one "最低票价" (G column) counts across
# the 展览, 演出 and 全部类型 sheets, matching refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 273
$ws.Range("F3").Value = 468
$ws.Range("F4").Value = 514
$ws.Range("F5").Value = 2414
$ws.Range("F7").Value = 76
$ws.Range("F9").Value = 1693
$ws.Range("F10").Value = 1693
$ws.Range("F16").Value = 856
$ws.Range("F20").Value = 7517
$ws.Range("F21").Value = 8504
$ws.Range("F24").Value = 420
$ws.Range("F42").Value = 1373
$ws.Range("F44").Value = 272
$ws.Range("F48").Value = 191
$ws.Range("F49").Value = 33

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 11
$ws.Range("G16").Value = 126

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 273
$ws.Range("F4").Value = 468
$ws.Range("F9").Value = 514
$ws.Range("F10").Value = 2414
$ws.Range("F12").Value = 76
$ws.Range("F14").Value = 1693
$ws.Range("F15").Value = 1693
$ws.Range("F19").Value = 856
$ws.Range("F25").Value = 7517
$ws.Range("F26").Value = 8504
$ws.Range("F42").Value = 1373
$ws.Range("F44").Value = 272
$ws.Range("F50").Value = 33

$wb.Save()
